# Daily auto-push data update.
# Two new timestamp rows are inserted right after the 2026/01/26 block
# (continuing that day, then starting 2026/01/27), shifting the existing
# 2026/12/29 .. 2027/01/05 block down by two rows. The sheet's used range
# grows from A1:D761 to A1:D763.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 720:721 - everything from the old row 720
# onward (2026/12/29 ...) shifts down to row 722 onward automatically.
$ws.Rows("720:721").Insert()

# New row 720: 2026/01/26 (Mon), 23:00, rank 201
# New row 721: 2026/01/27 (Tue), 01:00, rank 201
# Format the date column as text first so "2026/01/26" style strings are
# stored as literal text (matching the rest of column A) instead of being
# auto-converted to a date serial number, then clear the temporary format
# so the cells end up with no explicit style, like every other data row.
$ws.Range("A720:A721").NumberFormat = "@"

$ws.Range("A720").Value = "2026/01/26"
$ws.Range("B720").Value = "月"
$ws.Range("C720").Value = 23
$ws.Range("D720").Value = 201

$ws.Range("A721").Value = "2026/01/27"
$ws.Range("B721").Value = "火"
$ws.Range("C721").Value = 1
$ws.Range("D721").Value = 201

$ws.Range("A720:A721").ClearFormats()
